# "updates at end of 2 weeks"
#
# - Reposition the "ReasonforVisit" table group (Group 165) and the
#   "VisitType" table group (Group 168).
# - Remove the "DayofVisit" table group (Group 174) entirely.
# - Resize/move the "Elbow Connector 180" connector.
# - Resize/move the "Elbow Connector 183" connector to the values that used
#   to belong to "Elbow Connector 186", and remove "Elbow Connector 186".
#
# Note: PowerPoint's Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) while the underlying XML stores EMU integers. To avoid
# the COM layer's point->EMU rounding landing one EMU short, every target
# EMU value below is nudged by half an EMU (expressed in points) before
# being assigned.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

# --- Move "Group 165" (ReasonforVisit table): (5079332,609600) -> (5950875,397327)
$reasonForVisit = $s.Shapes.Item("Group 165")
$reasonForVisit.Left = EmuToPt 5950875
$reasonForVisit.Top  = EmuToPt 397327

# --- Move "Group 168" (VisitType table): (2835547,609600) -> (3663617,401180)
$visitType = $s.Shapes.Item("Group 168")
$visitType.Left = EmuToPt 3663617
$visitType.Top  = EmuToPt 401180

# --- Remove "Group 174" (DayofVisit table) entirely
$s.Shapes.Item("Group 174").Delete()

# --- Resize/move "Elbow Connector 180"
$conn180 = $s.Shapes.Item("Elbow Connector 180")
$conn180.Left   = EmuToPt 4948238
$conn180.Top    = EmuToPt 708711
$conn180.Width  = EmuToPt 679682
$conn180.Height = EmuToPt 1255695

# --- Resize/move "Elbow Connector 183" to the values formerly on
#     "Elbow Connector 186", then delete "Elbow Connector 186"
$conn183 = $s.Shapes.Item("Elbow Connector 183")
$conn183.Left   = EmuToPt 6272519
$conn183.Top    = EmuToPt 1003694
$conn183.Width  = EmuToPt 685801
$conn183.Height = EmuToPt 664142

$s.Shapes.Item("Elbow Connector 186").Delete()
